# ManageProducts.xlsx — "Recording flat and existing scripts update"
# Updates a handful of product rows on the "Input" sheet: new product
# names, adjusted carton dimensions, and a couple of Category2 swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Row 2 — NP-SC-SKU
$ws.Range("B2").Value = "prodbYcp"
$ws.Range("I2").Value = 70
$ws.Range("K2").Value = 11
$ws.Range("M2").Value = 12
$ws.Range("O2").Value = 13
$ws.Range("Q2").Value = "Safe/Vault"

# Row 3 — NP-MC-SKU
$ws.Range("B3").Value = "prodETnq"

# Row 5 — P-MC-SKU
$ws.Range("B5").Value = "prodVPcM"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 20
$ws.Range("P5").Value = 10
$ws.Range("R5").Value = "Other"

# Restore the view: scrolled right a bit, selection on P2
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("P2").Select()
